# Update "23 sectors" column (E) results for the Walktrap/correlations row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 was text "9 " (shared string) -> now a plain number 4
$ws.Range("E2").Value = 4

# E3 was number 17 -> now 5
$ws.Range("E3").Value = 5

# Move the active selection to E4, matching where the author's edit left off
$ws.Range("E4").Select()
